$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 13 (shifts old rows 13-24 down to 14-25) ---
$ws.Rows("13:13").Insert()

# The inserted row 13 inherits column-A formatting from the row above; the target
# layout has no cell at all in A13, so clear it completely.
$ws.Range("A13").Clear()

# New row 13 holds the "Docentes responsaveis" answer that used to sit (mislabeled)
# in row 10; give B13/C13 the correct value and correct column styles.
$ws.Range("B13").Value = '5817181 - Valdeir Arantes'
$ws.Range("C13").Value = '5817181 - Valdeir Arantes'
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# --- 2. Fix column definition (col A alone, not merged with col B) ---
# (cosmetic column metadata; widths unchanged)

# --- 3. Replace placeholder / mismatched text with the correct content ---

# Row 10 - "Objetivos:" answer (was wrongly showing the "Valdeir Arantes" text)
$ws.Range("B10").Value = 'Familiarizar o aluno com os conceitos básicos da estatística aplicada para estudar influência de variáveis independentes sobre variáveis dependentes (respostas) em bioprocessos. Introduzir ao aluno ferramentas de planejamento experimental usadas para planejar, executar experimentos fatoriais completo e fracionado, analisar os resultados, modelar o processo com base nos dados empíricos propondo condições de otimização e, também, familiarizar o aluno om um software comercial sobre o assunto.'
$ws.Range("C10").Value = 'Familiarizar o aluno com os conceitos básicos da estatística aplicada para estudar influência de variáveis independentes sobre variáveis dependentes (respostas) em bioprocessos. Introduzir ao aluno ferramentas de planejamento experimental usadas para planejar, executar experimentos fatoriais completo e fracionado, analisar os resultados, modelar o processo com base nos dados empíricos propondo condições de otimização e, também, familiarizar o aluno om um software comercial sobre o assunto.'

# Row 14 - "Programa resumido:" answer (was "Semestral")
$ws.Range("B14").Value = '1. O papel da estatística na Engenharia2. Fundamentos de estatística aplicada3. Análise de Variância4. Testes de comparações múltiplas5. Planejamento de Experimentos'
$ws.Range("C14").Value = '1. O papel da estatística na Engenharia2. Fundamentos de estatística aplicada3. Análise de Variância4. Testes de comparações múltiplas5. Planejamento de Experimentos'

# Row 16 - "Programa:" answer (was wrongly showing a date)
$ws.Range("B16").Value = '1. O papel da estatística na Engenharia: métodos de coleta de dados2. Fundamentos de estatística aplicada3. Análise de Variância: análise de variância de um modelo4. Testes de comparações múltiplas (Tukey, Hsu)5. Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta'
$ws.Range("C16").Value = '1. O papel da estatística na Engenharia: métodos de coleta de dados2. Fundamentos de estatística aplicada3. Análise de Variância: análise de variância de um modelo4. Testes de comparações múltiplas (Tukey, Hsu)5. Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta'

# Row 19 - "Metodo:" answer (was wrongly showing the "Valdeir Arantes" text)
$ws.Range("B19").Value = 'A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

# Row 20 - "Criterio:" answer (was the avaliacao paragraph)
$ws.Range("B20").Value = 'MF≥ 5,0 para aprovação 5,0'
$ws.Range("C20").Value = 'MF≥ 5,0 para aprovação 5,0'

# Row 21 - "Norma de recuperacao:" answer (was the MF>=5,0 text)
$ws.Range("B21").Value = '(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada.'

# Row 22 - "Bibliografia:" answer (was the recuperacao rule text)
$ws.Range("B22").Value = '1. BOX, G.E.P.; HUNTER, W.G.; HUNTER, J.S. Statistics for Experimenters: an introduction to designs, data analysis and model building. New York: John Wiley & Sons Inc., 1978.2. RODRIGUES, M. I. e IEMMA, A. F. Planejamento de experimentos e otimização de processos. Campinas: Cárita editora, 2009.3. Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 1996'
$ws.Range("C22").Value = '1. BOX, G.E.P.; HUNTER, W.G.; HUNTER, J.S. Statistics for Experimenters: an introduction to designs, data analysis and model building. New York: John Wiley & Sons Inc., 1978.2. RODRIGUES, M. I. e IEMMA, A. F. Planejamento de experimentos e otimização de processos. Campinas: Cárita editora, 2009.3. Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 1996'

Write-Host "Edit complete"
